$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Data" header moves from B1 to C1 (new column inserted to its left);
# copy B1's formatting over to C1 first, then move the value, then reset
# B1 back to a plain 0 (same style it already had).
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value2 = $ws.Range("B1").Value2
$ws.Range("B1").Value2 = 0

# Rows 2-4: the numeric data in column B moves over to column C.
$ws.Range("C2:C4").Value2 = $ws.Range("B2:B4").Value2
$ws.Range("B2:B4").ClearContents()

# Rows 5-7 are left untouched.
